$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 2-10 data (columns B-G) down to rows 3-11 (from bottom up to avoid overwrite issues)
for ($r = 10; $r -ge 2; $r--) {
    $src = $r
    $dst = $r + 1
    for ($c = 2; $c -le 7; $c++) {
        $ws.Cells.Item($dst, $c).Value = $ws.Cells.Item($src, $c).Value2
    }
}

# Set new values for row 2 (columns B-G)
$ws.Cells.Item(2, 2).Value = -0.08312210549351147
$ws.Cells.Item(2, 3).Value = 0.6473643844303238
$ws.Cells.Item(2, 4).Value = 0.6561623056018049
$ws.Cells.Item(2, 5).Value = 0.8100384593349904
$ws.Cells.Item(2, 6).Value = 0.8340433064168259
$ws.Cells.Item(2, 7).Value = 15
